$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.042.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.03%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.403.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.49%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'506.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.61%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'133.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.50%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.414.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.03%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +2.22%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -1.32%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.41%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -3.00%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.832.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.52%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'56.977.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.08%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'21.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.38%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +2.47%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.400.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.75%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'10.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.30%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +0.19%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'310.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.53%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.47%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.13%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'5.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.30%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'67.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +4.23%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -0.19%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.376"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.87%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.15%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'7.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.96%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'175.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.46%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0₃0728"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.05%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  -0.11%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +1.43%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.12%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +0.14%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.16%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'17.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.34%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.17%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.29%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.829"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.25%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'36.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.37%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +0.87%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'133.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.73%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +1.19%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'4.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.31%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.572"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.04%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Stellar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0915"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.55%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Bittensor"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'251.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.24%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +0.64%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +2.33%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'17.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +7.67%  "
$ws.Range("E51").Style = "Normal"

